# Add the next diagonal of values to the triangular error-series table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @{
    "K11" = 0.1683237681281231
    "J12" = 0.1722916656412322
    "I13" = 0.3079317558114735
    "H14" = 0.06712557395580883
    "G15" = 0.02179435870371246
    "F16" = -0.04506706323234141
    "E17" = -0.07465326558905801
    "D18" = -0.0928039223186989
    "C19" = -0.1108357465673982
    "B20" = -0.1624199859130616
}

foreach ($addr in $edits.Keys) {
    $ws.Range($addr).Value = $edits[$addr]
}
